# Update the monthly statistics table (StatBudget) with the real figures
# for MAI, JUIN and JUIL (previously placeholder zeros), and propagate the
# resulting running balance ("SOLD") into the remaining months of the year.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- MAI (table row 7) ---------------------------------------------------
$t.Cell(7, 3).Range.Text = "4138"
$t.Cell(7, 4).Range.Text = "3 724 200,00"
$t.Cell(7, 5).Range.Text = "238"
$t.Cell(7, 6).Range.Text = "214 200,00"
$t.Cell(7, 7).Range.Text = "4376"
$t.Cell(7, 8).Range.Text = "3 938 400,00"
$t.Cell(7, 9).Range.Text = "33 781 700,00"

# --- JUIN (table row 8) ---------------------------------------------------
$t.Cell(8, 3).Range.Text = "4292"
$t.Cell(8, 4).Range.Text = "4 292 000,00"
$t.Cell(8, 5).Range.Text = "286"
$t.Cell(8, 6).Range.Text = "286 000,00"
$t.Cell(8, 7).Range.Text = "4578"
$t.Cell(8, 8).Range.Text = "4 578 000,00"
$t.Cell(8, 9).Range.Text = "29 203 700,00"

# --- JUIL (table row 9) ---------------------------------------------------
$t.Cell(9, 3).Range.Text = "4306"
$t.Cell(9, 4).Range.Text = "4 306 000,00"
$t.Cell(9, 7).Range.Text = "4306"
$t.Cell(9, 8).Range.Text = "4 306 000,00"
$t.Cell(9, 9).Range.Text = "24 897 700,00"

# --- Remaining months (AOUT..DEC, table rows 10-14): only SOLD changes ---
for ($row = 10; $row -le 14; $row++) {
    $t.Cell($row, 9).Range.Text = "24 897 700,00"
}
